# The "Fecha" column (K) stored its dates as free text using a space as the
# separator (e.g. "09 08 24"). Re-format those values to use a slash
# separator instead ("09/08/24") for every data row, leaving the header
# (K1) untouched.
#
# The target text must stay a literal string (not get auto-converted into a
# real Excel date serial number), so the cells are pre-formatted as Text
# before the values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 11).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 23 }

$oldValue = "09 08 24"
$newValue = "09/08/24"

$dataRange = $ws.Range($ws.Cells.Item(2, 11), $ws.Cells.Item($lastRow, 11))
$dataRange.NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
